# Commit 57 Final Version 1.0.3.1
# Sample data error & model unique value change
#
# Fills in the previously-blank rows 41-48 (Sunday "D" sales) on Hoja1,
# mirroring the pattern used by the other day blocks (rows 2-9, 10-17, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the row-40 formatting (fonts/alignment per column) down onto the
# rows we are about to populate, matching the style used by every other
# day block in the sheet.
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G48").PasteSpecial(-4122)

# Día | Item | Price | Quantity | Cost | Categoria
$rows = @(
    @{ Row = 41; Dia = "D"; Item = "Primitiva";    Price = "177,00"; Qty = 177; Cost = 0; Cat = "Loteria" },
    @{ Row = 42; Dia = "D"; Item = "Joker Prim.";  Price = "7,00";   Qty = 7;   Cost = 0; Cat = "Loteria" },
    @{ Row = 43; Dia = "D"; Item = "Bono Loto";    Price = "133,00"; Qty = 266; Cost = 0; Cat = "Loteria" },
    @{ Row = 44; Dia = "D"; Item = "Gordo";        Price = "196,50"; Qty = 131; Cost = 0; Cat = "Loteria" },
    @{ Row = 45; Dia = "D"; Item = "Quiniela";     Price = "1,50";   Qty = 2;   Cost = 0; Cat = "Loteria" },
    @{ Row = 46; Dia = "D"; Item = "Quinigol";     Price = "1,00";   Qty = 1;   Cost = 0; Cat = "Loteria" },
    @{ Row = 47; Dia = "D"; Item = "Euromillones"; Price = "112,50"; Qty = 45;  Cost = 0; Cat = "Loteria" },
    @{ Row = 48; Dia = "D"; Item = "Lototurf";     Price = "3,00";   Qty = 3;   Cost = 0; Cat = "Loteria" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Dia
    $ws.Cells.Item($r.Row, 2).Value = $r.Item
    $ws.Cells.Item($r.Row, 3).Value = $r.Price
    $ws.Cells.Item($r.Row, 4).Value = $r.Qty
    $ws.Cells.Item($r.Row, 5).Value = $r.Cost
    $ws.Cells.Item($r.Row, 6).Value = $r.Cat
}

# Reflect the author's final view state: the newly filled-in rows
# (as whole rows) left selected.
$ws.Range("A41:XFD48").Select()
